$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = [double]"15.829186"
$ws.Range("H2").Value = [double]"47.487558"
$ws.Range("I2").Value = [double]"0.01520167221269649"
$ws.Range("J2").Value = [double]"0.01552195334947967"
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = [double]"197.2278263333334"
$ws.Range("N2").Value = [double]"591.683479"
$ws.Range("O2").Value = [double]"0.6783778564662776"
$ws.Range("P2").Value = [double]"0.6850369527608899"
$ws.Range("Q2").Value = [double]"3121.955947406032"
$ws.Range("R2").Value = [double]"28097.60352665428"
$ws.Range("S2").Value = [double]"0.01031247781035202"
$ws.Range("T2").Value = [double]"0.01063311162342424"

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = [double]"15.829186"
$ws.Range("H3").Value = [double]"47.487558"
$ws.Range("I3").Value = [double]"0.01520167221269649"
$ws.Range("J3").Value = [double]"0.01552195334947967"
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = [double]"0.896351"
$ws.Range("N3").Value = [double]"2.689053"
$ws.Range("O3").Value = [double]"0.003083057200020643"
$ws.Range("P3").Value = [double]"0.003113321122377543"
$ws.Range("Q3").Value = [double]"14.188506700286"
$ws.Range("R3").Value = [double]"127.696560302574"
$ws.Range("S3").Value = [double]"4.686762496770765e-05"
$ws.Range("T3").Value = [double]"4.83248252234939e-05"

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = [double]"15.829186"
$ws.Range("H4").Value = [double]"47.487558"
$ws.Range("I4").Value = [double]"0.01520167221269649"
$ws.Range("J4").Value = [double]"0.01552195334947967"
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = [double]"42.82536200000001"
$ws.Range("N4").Value = [double]"128.476086"
$ws.Range("O4").Value = [double]"0.1473006006102414"
$ws.Range("P4").Value = [double]"0.1487465335432934"
$ws.Range("Q4").Value = [double]"677.8906206153321"
$ws.Range("R4").Value = [double]"6101.015585537988"
$ws.Range("S4").Value = [double]"0.00223921544721021"
$ws.Range("T4").Value = [double]"0.002308836754555813"

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = [double]"15.829186"
$ws.Range("H5").Value = [double]"47.487558"
$ws.Range("I5").Value = [double]"0.01520167221269649"
$ws.Range("J5").Value = [double]"0.01552195334947967"
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = [double]"41.30642933333333"
$ws.Range("N5").Value = [double]"123.919288"
$ws.Range("O5").Value = [double]"0.1420761335272424"
$ws.Range("P5").Value = [double]"0.1434707820189434"
$ws.Range("Q5").Value = [double]"653.8471529131893"
$ws.Range("R5").Value = [double]"5884.624376218704"
$ws.Range("S5").Value = [double]"0.002159794811128436"
$ws.Range("T5").Value = [double]"0.002226946785511407"

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = [double]"15.829186"
$ws.Range("H6").Value = [double]"47.487558"
$ws.Range("I6").Value = [double]"0.01520167221269649"
$ws.Range("J6").Value = [double]"0.01552195334947967"
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = [double]"8.478501"
$ws.Range("N6").Value = [double]"16.957002"
$ws.Range("O6").Value = [double]"0.02916235219621802"
$ws.Range("P6").Value = [double]"0.01963241055449567"
$ws.Range("Q6").Value = [double]"134.207769330186"
$ws.Range("R6").Value = [double]"805.246615981116"
$ws.Range("S6").Value = [double]"0.0004433165190381159"
$ws.Range("T6").Value = [double]"0.0003047333607647141"

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = [double]"155.500389"
$ws.Range("H7").Value = [double]"466.501167"
$ws.Range("I7").Value = [double]"0.1493359129474374"
$ws.Range("J7").Value = [double]"0.1524822428572096"
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = [double]"197.2278263333334"
$ws.Range("N7").Value = [double]"591.683479"
$ws.Range("O7").Value = [double]"0.6783778564662776"
$ws.Range("P7").Value = [double]"0.6850369527608899"
$ws.Range("Q7").Value = [double]"30669.00371645778"
$ws.Range("R7").Value = [double]"276021.03344812"
$ws.Range("S7").Value = [double]"0.1013061765187172"
$ws.Range("T7").Value = [double]"0.1044559709970488"

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = [double]"155.500389"
$ws.Range("H8").Value = [double]"466.501167"
$ws.Range("I8").Value = [double]"0.1493359129474374"
$ws.Range("J8").Value = [double]"0.1524822428572096"
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = [double]"0.896351"
$ws.Range("N8").Value = [double]"2.689053"
$ws.Range("O8").Value = [double]"0.003083057200020643"
$ws.Range("P8").Value = [double]"0.003113321122377543"
$ws.Range("Q8").Value = [double]"139.382929180539"
$ws.Range("R8").Value = [double]"1254.446362624851"
$ws.Range("S8").Value = [double]"0.0004604111616342529"
$ws.Range("T8").Value = [double]"0.0004747261874748527"

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = [double]"155.500389"
$ws.Range("H9").Value = [double]"466.501167"
$ws.Range("I9").Value = [double]"0.1493359129474374"
$ws.Range("J9").Value = [double]"0.1524822428572096"
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = [double]"42.82536200000001"
$ws.Range("N9").Value = [double]"128.476086"
$ws.Range("O9").Value = [double]"0.1473006006102414"
$ws.Range("P9").Value = [double]"0.1487465335432934"
$ws.Range("Q9").Value = [double]"6659.360450065819"
$ws.Range("R9").Value = [double]"59934.24405059237"
$ws.Range("S9").Value = [double]"0.02199726966983625"
$ws.Range("T9").Value = [double]"0.02268120505191653"

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = [double]"155.500389"
$ws.Range("H10").Value = [double]"466.501167"
$ws.Range("I10").Value = [double]"0.1493359129474374"
$ws.Range("J10").Value = [double]"0.1524822428572096"
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = [double]"41.30642933333333"
$ws.Range("N10").Value = [double]"123.919288"
$ws.Range("O10").Value = [double]"0.1420761335272424"
$ws.Range("P10").Value = [double]"0.1434707820189434"
$ws.Range("Q10").Value = [double]"6423.165829534345"
$ws.Range("R10").Value = [double]"57808.4924658091"
$ws.Range("S10").Value = [double]"0.02121706910833276"
$ws.Range("T10").Value = [double]"0.02187674662672631"

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = [double]"155.500389"
$ws.Range("H11").Value = [double]"466.501167"
$ws.Range("I11").Value = [double]"0.1493359129474374"
$ws.Range("J11").Value = [double]"0.1524822428572096"
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = [double]"8.478501"
$ws.Range("N11").Value = [double]"16.957002"
$ws.Range("O11").Value = [double]"0.02916235219621802"
$ws.Range("P11").Value = [double]"0.01963241055449567"
$ws.Range("Q11").Value = [double]"1318.410203636889"
$ws.Range("R11").Value = [double]"7910.461221821334"
$ws.Range("S11").Value = [double]"0.004354986488916925"
$ws.Range("T11").Value = [double]"0.002993593994043053"

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = [double]"179.1193723333333"
$ws.Range("H12").Value = [double]"537.358117"
$ws.Range("I12").Value = [double]"0.1720185728536685"
$ws.Range("J12").Value = [double]"0.1756427994052303"
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = [double]"197.2278263333334"
$ws.Range("N12").Value = [double]"591.683479"
$ws.Range("O12").Value = [double]"0.6783778564662776"
$ws.Range("P12").Value = [double]"0.6850369527608899"
$ws.Range("Q12").Value = [double]"35327.32445949434"
$ws.Range("R12").Value = [double]"317945.920135449"
$ws.Range("S12").Value = [double]"0.1166935907248599"
$ws.Range("T12").Value = [double]"0.1203218080789512"

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = [double]"179.1193723333333"
$ws.Range("H13").Value = [double]"537.358117"
$ws.Range("I13").Value = [double]"0.1720185728536685"
$ws.Range("J13").Value = [double]"0.1756427994052303"
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = [double]"0.896351"
$ws.Range("N13").Value = [double]"2.689053"
$ws.Range("O13").Value = [double]"0.003083057200020643"
$ws.Range("P13").Value = [double]"0.003113321122377543"
$ws.Range("Q13").Value = [double]"160.5538285103557"
$ws.Range("R13").Value = [double]"1444.984456593201"
$ws.Range("S13").Value = [double]"0.0005303430995737783"
$ws.Range("T13").Value = [double]"0.0005468324373818251"

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = [double]"179.1193723333333"
$ws.Range("H14").Value = [double]"537.358117"
$ws.Range("I14").Value = [double]"0.1720185728536685"
$ws.Range("J14").Value = [double]"0.1756427994052303"
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = [double]"42.82536200000001"
$ws.Range("N14").Value = [double]"128.476086"
$ws.Range("O14").Value = [double]"0.1473006006102414"
$ws.Range("P14").Value = [double]"0.1487465335432934"
$ws.Range("Q14").Value = [double]"7670.851961387785"
$ws.Range("R14").Value = [double]"69037.66765249007"
$ws.Range("S14").Value = [double]"0.02533843909746194"
$ws.Range("T14").Value = [double]"0.02612625755336804"

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = [double]"179.1193723333333"
$ws.Range("H15").Value = [double]"537.358117"
$ws.Range("I15").Value = [double]"0.1720185728536685"
$ws.Range("J15").Value = [double]"0.1756427994052303"
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = [double]"41.30642933333333"
$ws.Range("N15").Value = [double]"123.919288"
$ws.Range("O15").Value = [double]"0.1420761335272424"
$ws.Range("P15").Value = [double]"0.1434707820189434"
$ws.Range("Q15").Value = [double]"7398.781695517855"
$ws.Range("R15").Value = [double]"66589.0352596607"
$ws.Range("S15").Value = [double]"0.02443973372592348"
$ws.Range("T15").Value = [double]"0.0251996097866648"

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = [double]"179.1193723333333"
$ws.Range("H16").Value = [double]"537.358117"
$ws.Range("I16").Value = [double]"0.1720185728536685"
$ws.Range("J16").Value = [double]"0.1756427994052303"
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = [double]"8.478501"
$ws.Range("N16").Value = [double]"16.957002"
$ws.Range("O16").Value = [double]"0.02916235219621802"
$ws.Range("P16").Value = [double]"0.01963241055449567"
$ws.Range("Q16").Value = [double]"1518.663777447539"
$ws.Range("R16").Value = [double]"9111.982664685234"
$ws.Range("S16").Value = [double]"0.00501646620584947"
$ws.Range("T16").Value = [double]"0.003448291548864408"

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = [double]"626.3728126666666"
$ws.Range("H17").Value = [double]"1879.118438"
$ws.Range("I17").Value = [double]"0.6015416194555684"
$ws.Range("J17").Value = [double]"0.6142153852759307"
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = [double]"197.2278263333334"
$ws.Range("N17").Value = [double]"591.683479"
$ws.Range("O17").Value = [double]"0.6783778564662776"
$ws.Range("P17").Value = [double]"0.6850369527608899"
$ws.Range("Q17").Value = [double]"123538.1483165429"
$ws.Range("R17").Value = [double]"1111843.334848886"
$ws.Range("S17").Value = [double]"0.4080725143815218"
$ws.Range("T17").Value = [double]"0.4207602358682795"

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = [double]"626.3728126666666"
$ws.Range("H18").Value = [double]"1879.118438"
$ws.Range("I18").Value = [double]"0.6015416194555684"
$ws.Range("J18").Value = [double]"0.6142153852759307"
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = [double]"0.896351"
$ws.Range("N18").Value = [double]"2.689053"
$ws.Range("O18").Value = [double]"0.003083057200020643"
$ws.Range("P18").Value = [double]"0.003113321122377543"
$ws.Range("Q18").Value = [double]"561.4498970065792"
$ws.Range("R18").Value = [double]"5053.049073059214"
$ws.Range("S18").Value = [double]"0.001854587220974568"
$ws.Range("T18").Value = [double]"0.001912249732668815"

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = [double]"626.3728126666666"
$ws.Range("H19").Value = [double]"1879.118438"
$ws.Range("I19").Value = [double]"0.6015416194555684"
$ws.Range("J19").Value = [double]"0.6142153852759307"
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = [double]"42.82536200000001"
$ws.Range("N19").Value = [double]"128.476086"
$ws.Range("O19").Value = [double]"0.1473006006102414"
$ws.Range("P19").Value = [double]"0.1487465335432934"
$ws.Range("Q19").Value = [double]"26824.64244940819"
$ws.Range("R19").Value = [double]"241421.7820446737"
$ws.Range("S19").Value = [double]"0.08860744183786247"
$ws.Range("T19").Value = [double]"0.09136240940875309"

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = [double]"626.3728126666666"
$ws.Range("H20").Value = [double]"1879.118438"
$ws.Range("I20").Value = [double]"0.6015416194555684"
$ws.Range("J20").Value = [double]"0.6142153852759307"
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = [double]"41.30642933333333"
$ws.Range("N20").Value = [double]"123.919288"
$ws.Range("O20").Value = [double]"0.1420761335272424"
$ws.Range("P20").Value = [double]"0.1434707820189434"
$ws.Range("Q20").Value = [double]"25873.2243227369"
$ws.Range("R20").Value = [double]"232859.0189046321"
$ws.Range("S20").Value = [double]"0.08546470744796295"
$ws.Range("T20").Value = [double]"0.08812196165360442"

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = [double]"626.3728126666666"
$ws.Range("H21").Value = [double]"1879.118438"
$ws.Range("I21").Value = [double]"0.6015416194555684"
$ws.Range("J21").Value = [double]"0.6142153852759307"
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = [double]"8.478501"
$ws.Range("N21").Value = [double]"16.957002"
$ws.Range("O21").Value = [double]"0.02916235219621802"
$ws.Range("P21").Value = [double]"0.01963241055449567"
$ws.Range("Q21").Value = [double]"5310.702518567145"
$ws.Range("R21").Value = [double]"31864.21511140287"
$ws.Range("S21").Value = [double]"0.01754236856724664"
$ws.Range("T21").Value = [double]"0.0120585286126248"

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = [double]"64.4575005"
$ws.Range("H22").Value = [double]"128.915001"
$ws.Range("I22").Value = [double]"0.06190222253062919"
$ws.Range("J22").Value = [double]"0.04213761911214986"
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = [double]"197.2278263333334"
$ws.Range("N22").Value = [double]"591.683479"
$ws.Range("O22").Value = [double]"0.6783778564662776"
$ws.Range("P22").Value = [double]"0.6850369527608899"
$ws.Range("Q22").Value = [double]"12712.81271449475"
$ws.Range("R22").Value = [double]"76276.87628696847"
$ws.Range("S22").Value = [double]"0.04199309703082674"
$ws.Range("T22").Value = [double]"0.02886582619318618"

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = [double]"64.4575005"
$ws.Range("H23").Value = [double]"128.915001"
$ws.Range("I23").Value = [double]"0.06190222253062919"
$ws.Range("J23").Value = [double]"0.04213761911214986"
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = [double]"0.896351"
$ws.Range("N23").Value = [double]"2.689053"
$ws.Range("O23").Value = [double]"0.003083057200020643"
$ws.Range("P23").Value = [double]"0.003113321122377543"
$ws.Range("Q23").Value = [double]"57.77654503067549"
$ws.Range("R23").Value = [double]"346.659270184053"
$ws.Range("S23").Value = [double]"0.0001908480928703364"
$ws.Range("T23").Value = [double]"0.0001311879396285558"

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = [double]"64.4575005"
$ws.Range("H24").Value = [double]"128.915001"
$ws.Range("I24").Value = [double]"0.06190222253062919"
$ws.Range("J24").Value = [double]"0.04213761911214986"
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = [double]"42.82536200000001"
$ws.Range("N24").Value = [double]"128.476086"
$ws.Range("O24").Value = [double]"0.1473006006102414"
$ws.Range("P24").Value = [double]"0.1487465335432934"
$ws.Range("Q24").Value = [double]"2760.415792527681"
$ws.Range("R24").Value = [double]"16562.49475516608"
$ws.Range("S24").Value = [double]"0.009118234557870494"
$ws.Range("T24").Value = [double]"0.00626782477469992"

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = [double]"64.4575005"
$ws.Range("H25").Value = [double]"128.915001"
$ws.Range("I25").Value = [double]"0.06190222253062919"
$ws.Range("J25").Value = [double]"0.04213761911214986"
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = [double]"41.30642933333333"
$ws.Range("N25").Value = [double]"123.919288"
$ws.Range("O25").Value = [double]"0.1420761335272424"
$ws.Range("P25").Value = [double]"0.1434707820189434"
$ws.Range("Q25").Value = [double]"2662.509189406548"
$ws.Range("R25").Value = [double]"15975.05513643929"
$ws.Range("S25").Value = [double]"0.008794828433894742"
$ws.Range("T25").Value = [double]"0.006045517166436518"

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = [double]"64.4575005"
$ws.Range("H26").Value = [double]"128.915001"
$ws.Range("I26").Value = [double]"0.06190222253062919"
$ws.Range("J26").Value = [double]"0.04213761911214986"
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = [double]"8.478501"
$ws.Range("N26").Value = [double]"16.957002"
$ws.Range("O26").Value = [double]"0.02916235219621802"
$ws.Range("P26").Value = [double]"0.01963241055449567"
$ws.Range("Q26").Value = [double]"546.5029824467505"
$ws.Range("R26").Value = [double]"2186.011929787002"
$ws.Range("S26").Value = [double]"0.001805214415166871"
$ws.Range("T26").Value = [double]"0.0008272630381986893"
